# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.933.20"
$ws.Range("E2").Value = "  +0.78%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.879.90"
$ws.Range("E3").Value = "  +0.33%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.019"
$ws.Range("E4").Value = "  +1.46%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.20"
$ws.Range("E5").Value = "  +0.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.017"
$ws.Range("E6").Value = "  +1.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4672"
$ws.Range("E7").Value = "  -1.48%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3909"
$ws.Range("E8").Value = "  -1.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.78"
$ws.Range("E9").Value = "  -1.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07936"
$ws.Range("E10").Value = "  -1.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.004"
$ws.Range("E11").Value = "  -1.67%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.56"
$ws.Range("E12").Value = "  -1.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.901.94"
$ws.Range("E13").Value = "  +0.97%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.941"
$ws.Range("E14").Value = "  -0.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.086"
$ws.Range("E15").Value = "  -0.76%  "

$ws.Range("E16").Value = "  +1.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06797"
$ws.Range("E17").Value = "  +2.52%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.44"
$ws.Range("E18").Value = "  +0.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.98"
$ws.Range("E20").Value = "  -1.67%  "

$ws.Range("E21").Value = "  +1.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.936.18"
$ws.Range("E22").Value = "  +0.72%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.461"
$ws.Range("E23").Value = "  -0.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.95"
$ws.Range("E24").Value = "  -0.60%  "

$ws.Range("E25").Value = "  +2.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.111.93"
$ws.Range("E26").Value = "  +0.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.43"
$ws.Range("E27").Value = "  +1.96%  "

$ws.Range("E28").Value = "  -1.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.066"
$ws.Range("E29").Value = "  -1.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.445"
$ws.Range("E30").Value = "  -2.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.68"
$ws.Range("E31").Value = "  -1.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09522"
$ws.Range("E32").Value = "  -0.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9530"
$ws.Range("E33").Value = "  -1.60%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.662"
$ws.Range("E34").Value = "  +0.83%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.312"
$ws.Range("E35").Value = "  +0.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.344"
$ws.Range("E36").Value = "  -7.35%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06111"
$ws.Range("E37").Value = "  +0.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02233"
$ws.Range("E38").Value = "  -1.17%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.199"
$ws.Range("E39").Value = "  -1.98%  "

$ws.Range("E40").Value = "  +1.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.084"
$ws.Range("E41").Value = "  -0.90%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5867"
$ws.Range("E42").Value = "  -1.93%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1892"
$ws.Range("E43").Value = "  -0.84%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.13"
$ws.Range("E44").Value = "  -1.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.275"
$ws.Range("E45").Value = "  +1.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5622"
$ws.Range("E46").Value = "  -1.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.09"
$ws.Range("E47").Value = "  -1.69%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.397"
$ws.Range("E48").Value = "  -0.26%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.911"
$ws.Range("E49").Value = "  -1.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06861"
$ws.Range("E50").Value = "  +0.61%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "113.39"
$ws.Range("E51").Value = "  +0.88%  "
